# Update for release to deploy 0.1.1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B): 0.1.0 -> 0.1.1
$ws.Range("B3").Value = "0.1.1"

# Update Date value (row 8, column B)
$ws.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Insert a new row for "Jurisdiction" after the "Contact" row (row 10),
# pushing all subsequent rows down by one.
$ws.Rows.Item(11).Insert()

# Copy the formatting from the row above (Contact row) onto the new row
# so the new cells carry the same style as the rest of the table.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Populate the new Jurisdiction row: Property = "Jurisdiction", Value = "" (empty)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
